$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 5 de Mayo de 2020 a las 07:03"

# Row 63 - Hungria
$ws.Range("B63").Value = 3065
$ws.Range("C63").Value = 30
$ws.Range("D63").Value = 709
$ws.Range("E63").Value = 1993
$ws.Range("G63").Value = 12
$ws.Range("H63").Value = 363

# Row 64 - Tailandia
$ws.Range("B64").Value = 2988
$ws.Range("C64").Value = 1
$ws.Range("D64").Value = 2747
$ws.Range("E64").Value = 187

# Row 96 - Kirguistan
$ws.Range("B96").Value = 843
$ws.Range("C96").Value = 13
$ws.Range("D96").Value = 600
$ws.Range("E96").Value = 232
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 11

# Row 115 - Maldivas
$ws.Range("B115").Value = 551
$ws.Range("C115").Value = 10
$ws.Range("E115").Value = 532

# Row 175 - Mongolia
$ws.Range("D175").Value = 13
$ws.Range("E175").Value = 28
